# Update the "取得日時" (retrieved at) timestamp column on the ランサーズ sheet
# for all existing data rows (A2:A11) from "2025-10-27 01:24:39" to
# "2025-10-27 01:55:00" (commit: "Append: 2025-10-27 01:55 JST").

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ランサーズ")

$newTimestamp = "2025-10-27 01:55:00"

$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row
for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 1)
    if ($cell.Value -ne $null -and $cell.Value -ne "") {
        $cell.Value = $newTimestamp
    }
}
